# Refresh the cryptos table (rows 2-51) with the latest scraped
# price/volume snapshot. Coin/Link (columns B/C) only change where
# the underlying ranking reshuffled; Price/Volume (D/E) are updated
# for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.381.91"
$ws.Cells.Item(2, 5).Value = "  -0.51%  "

$ws.Cells.Item(3, 4).Value = "1.721.94"
$ws.Cells.Item(3, 5).Value = "  -0.38%  "

$ws.Cells.Item(4, 5).Value = "  +0.02%  "

$ws.Cells.Item(5, 4).Value = "'244.18"
$ws.Cells.Item(5, 5).Value = "  -0.21%  "

$ws.Cells.Item(6, 4).Value = "'1.000"
$ws.Cells.Item(6, 5).Value = "  -0.04%  "

$ws.Cells.Item(7, 4).Value = "'0.4877"
$ws.Cells.Item(7, 5).Value = "  +1.84%  "

$ws.Cells.Item(8, 5).Value = "  -2.24%  "

$ws.Cells.Item(9, 4).Value = "'0.06203"
$ws.Cells.Item(9, 5).Value = "  -0.34%  "

$ws.Cells.Item(10, 4).Value = "1.725.39"
$ws.Cells.Item(10, 5).Value = "  -0.21%  "

$ws.Cells.Item(11, 4).Value = "'0.07017"
$ws.Cells.Item(11, 5).Value = "  -2.12%  "

$ws.Cells.Item(12, 4).Value = "'15.46"
$ws.Cells.Item(12, 5).Value = "  -1.34%  "

$ws.Cells.Item(13, 4).Value = "'4.532"
$ws.Cells.Item(13, 5).Value = "  +0.07%  "

$ws.Cells.Item(14, 4).Value = "'0.5960"
$ws.Cells.Item(14, 5).Value = "  -3.03%  "

$ws.Cells.Item(15, 4).Value = "'77.16"
$ws.Cells.Item(15, 5).Value = "  +0.34%  "

$ws.Cells.Item(16, 5).Value = "  +0.01%  "

$ws.Cells.Item(17, 4).Value = "26.396.44"
$ws.Cells.Item(17, 5).Value = "  -0.46%  "

$ws.Cells.Item(18, 5).Value = "  +0.05%  "

$ws.Cells.Item(19, 4).Value = "'0.000007213"
$ws.Cells.Item(19, 5).Value = "  +3.55%  "

$ws.Cells.Item(20, 4).Value = "'11.36"
$ws.Cells.Item(20, 5).Value = "  -2.35%  "

$ws.Cells.Item(21, 4).Value = "1.949.15"
$ws.Cells.Item(21, 5).Value = "  -0.15%  "

$ws.Cells.Item(22, 4).Value = "'4.476"
$ws.Cells.Item(22, 5).Value = "  -1.07%  "

$ws.Cells.Item(23, 4).Value = "'8.554"
$ws.Cells.Item(23, 5).Value = "  -4.04%  "

$ws.Cells.Item(24, 4).Value = "'5.166"
$ws.Cells.Item(24, 5).Value = "  -2.17%  "

$ws.Cells.Item(25, 4).Value = "'137.37"
$ws.Cells.Item(25, 5).Value = "  +0.67%  "

$ws.Cells.Item(26, 4).Value = "'15.22"
$ws.Cells.Item(26, 5).Value = "  -0.84%  "

$ws.Cells.Item(27, 4).Value = "'1.417"
$ws.Cells.Item(27, 5).Value = "  +0.94%  "

$ws.Cells.Item(28, 5).Value = "  +0.63%  "

$ws.Cells.Item(29, 5).Value = "  -4.35%  "

$ws.Cells.Item(30, 4).Value = "'3.959"
$ws.Cells.Item(30, 5).Value = "  -0.58%  "

$ws.Cells.Item(31, 4).Value = "'0.07942"
$ws.Cells.Item(31, 5).Value = "  -0.30%  "

$ws.Cells.Item(32, 4).Value = "'3.679"
$ws.Cells.Item(32, 5).Value = "  -0.83%  "

$ws.Cells.Item(33, 4).Value = "'0.04528"
$ws.Cells.Item(33, 5).Value = "  -1.32%  "

$ws.Cells.Item(34, 2).Value = "HuobiToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(34, 4).Value = "'2.615"
$ws.Cells.Item(34, 5).Value = "  -0.14%  "

$ws.Cells.Item(35, 2).Value = "ARBITRUM"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(35, 4).Value = "'0.9939"
$ws.Cells.Item(35, 5).Value = "  -0.01%  "

$ws.Cells.Item(36, 2).Value = "ImmutableX"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(36, 4).Value = "'0.6208"
$ws.Cells.Item(36, 5).Value = "  -1.59%  "

$ws.Cells.Item(37, 2).Value = "TrustWalletToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(37, 4).Value = "'0.9066"
$ws.Cells.Item(37, 5).Value = "  -0.78%  "

$ws.Cells.Item(38, 2).Value = "RenderToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(38, 4).Value = "'1.978"
$ws.Cells.Item(38, 5).Value = "  -5.20%  "

$ws.Cells.Item(39, 2).Value = "MXToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(39, 4).Value = "'2.395"
$ws.Cells.Item(39, 5).Value = "  -0.46%  "

$ws.Cells.Item(40, 2).Value = "PaxDollar"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(40, 4).Value = "'0.9999"
$ws.Cells.Item(40, 5).Value = "  -0.31%  "

$ws.Cells.Item(41, 2).Value = "VeChain"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(41, 4).Value = "'0.01485"
$ws.Cells.Item(41, 5).Value = "  -1.09%  "

$ws.Cells.Item(42, 2).Value = "Quant"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(42, 4).Value = "'100.25"
$ws.Cells.Item(42, 5).Value = "  -4.21%  "

$ws.Cells.Item(43, 2).Value = "FraxShare"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(43, 4).Value = "'5.386"
$ws.Cells.Item(43, 5).Value = "  -3.49%  "

$ws.Cells.Item(44, 2).Value = "TheSandbox"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(44, 4).Value = "'0.3844"
$ws.Cells.Item(44, 5).Value = "  -0.86%  "

$ws.Cells.Item(45, 2).Value = "Aptos"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(45, 4).Value = "'6.725"
$ws.Cells.Item(45, 5).Value = "  -3.65%  "

$ws.Cells.Item(46, 2).Value = "Algorand"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(46, 4).Value = "'0.1151"
$ws.Cells.Item(46, 5).Value = "  -2.66%  "

$ws.Cells.Item(47, 2).Value = "Cronos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(47, 4).Value = "'0.05354"
$ws.Cells.Item(47, 5).Value = "  +0.14%  "

$ws.Cells.Item(48, 2).Value = "Elrond"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(48, 4).Value = "'30.08"
$ws.Cells.Item(48, 5).Value = "  -2.97%  "

$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).Value = "'7.681"
$ws.Cells.Item(49, 5).Value = "  -2.19%  "

$ws.Cells.Item(50, 2).Value = "NEARProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(50, 4).Value = "'1.241"
$ws.Cells.Item(50, 5).Value = "  -1.42%  "

$ws.Cells.Item(51, 2).Value = "Aave"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(51, 4).Value = "'51.02"
$ws.Cells.Item(51, 5).Value = "  -0.28%  "
